# ---------------------------------------------------------------------------
# "modelo de diagnostico y seguimiento alimentario"
#
# Inserts a new block of paragraphs right after the
# "Revaloración (proceso cíclico)." bullet (end of the "Fase de inicio"
# list) and before "Fase de construcción.":
#
#   - an intro paragraph "Las herramientas utilizadas incluyen:"
#   - four new bulleted sub-items (same bullet list / numId as the other
#     sub-bullets already in the document)
#
# The trailing "_GoBack" bookmark that used to sit on the "Pruebas."
# paragraph is moved to sit on the new last paragraph of the inserted
# block (matching Word's usual behaviour of leaving _GoBack at the most
# recent edit point).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Insert the new paragraphs right after "Revaloración (proceso cíclico)." ---

$findRange = $d.Content
$findRange.Find.Execute("Revaloración (proceso cíclico).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $findRange.Paragraphs(1)
$anchorRange = $anchorPara.Range

# Collapse to just before the paragraph mark so the new content is inserted
# *inside* a fresh paragraph following the anchor, without merging into the
# next paragraph's properties.
$insertPoint = $d.Range($anchorRange.End - 1, $anchorRange.End - 1)

$newParagraphsXml = @'
<w:p>
  <w:pPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:ind w:left="708" w:firstLine="708"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Las herramientas uti</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>lizadas incluyen:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Esquema del plato del buen comer</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Esquema de la jarra del bien beber.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Herramienta para cálculos dietéticos.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Tabla de equivalencias de alimentos del Sistema de equivalentes mexicanos.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$insertPoint.InsertXML($newParagraphsXml)

# --- 2. Drop the old "_GoBack" bookmark that used to sit on "Pruebas." ---
# (it has now moved to the new last paragraph inserted above, so the old
# one must be removed to avoid a duplicate bookmark name/id).

$pruebasRange = $d.Content
$pruebasRange.Find.Execute("Pruebas.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pruebasPara = $pruebasRange.Paragraphs(1)
$pruebasParaRange = $pruebasPara.Range

$pruebasXml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="6"/>
    </w:numPr>
    <w:spacing w:before="120" w:after="120" w:line="360" w:lineRule="auto"/>
    <w:contextualSpacing w:val="0"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>Pruebas.</w:t>
  </w:r>
</w:p>
'@

$pruebasParaRange.InsertXML($pruebasXml)

Write-Output "Edit applied."
